# Croatia HNL - base update for 07-04-2024 22:30 run
# Inserts two new earlier fixtures (rows 142-143) and re-computes the closing
# lines for the two fixtures that used to occupy those rows (now shifted to 144-145).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 144 and 145 are brand new (the sheet previously ended at row 143).
# Seed their column-A / column-E formatting (bold+border id, date display) from the
# existing rows 142/143 via a formats-only paste so no new style entries are minted
# (styles.xml stays byte-identical to before).
$ws.Cells.Item(142, 1).Copy()
$ws.Cells.Item(144, 1).PasteSpecial(-4122)
$ws.Cells.Item(145, 1).PasteSpecial(-4122)
$ws.Cells.Item(142, 5).Copy()
$ws.Cells.Item(144, 5).PasteSpecial(-4122)
$ws.Cells.Item(145, 5).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Full target row contents, column A (id) through column AC (PL_AhUnder).
$rowsData = @(
    @{ Row = 142; Vals = @(140, 6787897, "Croatia HNL", "Croatia HNL", 45388.49305555555, "NK Lokomotiva Zagreb", "HNK Gorica", 1, 1, "D", 1.615, 3.5, 4.75, 1.727, 3.3, 4.2, -0.5, 1.825, 2.025, 2.25, 1.975, 1.875, -1, 2.3, -1, -1, 1.025, -0.5, 0.4375) }
    @{ Row = 143; Vals = @(141, 6788939, "Croatia HNL", "Croatia HNL", 45388.58333333334, "Slaven Belupo", "NK Varazdin", 0, 1, "A", 2.25, 3.25, 2.75, 2.15, 3.25, 2.9, -0.25, 2, 1.85, 2.5, 2.05, 1.8, -1, -1, 1.9, -1, 0.8500000000000001, -1, 0.8) }
    @{ Row = 144; Vals = @(142, 6788941, "Croatia HNL", "Croatia HNL", 45389.5, "Dinamo Zagreb", "Istra 1961", 4, 1, "H", 1.2, 6.5, 8, 1.166, 7.5, 10, -2, 1.975, 1.875, 3, 1.975, 1.875, 0.1659999999999999, -1, -1, 0.9750000000000001, -1, 0.9750000000000001, -1) }
    @{ Row = 145; Vals = @(143, 6788940, "Croatia HNL", "Croatia HNL", 45389.60416666666, "HNK Rijeka", "Hajduk Split", 1, 0, "H", 2.1, 3.25, 3, 1.8, 3.4, 3.75, -0.5, 1.925, 1.925, 2.25, 1.925, 1.925, 0.8, -1, -1, 0.925, -1, -1, 0.925) }
)

foreach ($entry in $rowsData) {
    $col = 1
    foreach ($v in $entry.Vals) {
        $ws.Cells.Item($entry.Row, $col).Value = $v
        $col = $col + 1
    }
}

